$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing address line "919 Story Road, San Jose CA 95122" into
#    two separate paragraphs: "919 Story Road" and "San Jose, CA 95122".
#    (The same address text also appears later inside a table and must stay
#    untouched, so target the addressee paragraph directly rather than using
#    a document-wide Find.)
$addrPara = $d.Paragraphs.Item(7)
if ($addrPara.Range.Text -like "919 Story Road, San Jose CA 95122*") {
    $addrPara.Range.Text = "919 Story Road`r"
    $newPara = $d.Paragraphs.Item(8)
    $newPara.Range.Text = "San Jose, CA 95122"
}

# 3. Remove the blank "NoSpacing" paragraph that used to follow the
#    "Board of Directors" signature line.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Board of Directors*") {
        $blankPara = $d.Paragraphs.Item($i + 1)
        $blankPara.Range.Delete()
        break
    }
}
